# Insert a new weekly price record as row 6 of the Albahaca sheet.
# All existing records from row 6 downward shift down by one row
# (handled automatically by EntireRow Insert), so the rest of the
# table is preserved unchanged, just moved down one position, and
# the sheet dimension grows from A1:R84 to A1:R85.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 6, pushing rows 6..84 down to 7..85.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new data record.
$ws.Cells.Item(6, 1).Value2 = 8
$ws.Cells.Item(6, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(6, 3).Value2 = "Coquimbo"
$ws.Cells.Item(6, 4).Value2 = 44685
$ws.Cells.Item(6, 5).Value2 = 4
$ws.Cells.Item(6, 6).Value2 = 100112052
$ws.Cells.Item(6, 7).Value2 = "Albahaca"
$ws.Cells.Item(6, 8).Value2 = "Sin especificar"
$ws.Cells.Item(6, 9).Value2 = "Primera"
$ws.Cells.Item(6, 10).Value2 = 2000
$ws.Cells.Item(6, 11).Value2 = 5000
$ws.Cells.Item(6, 12).Value2 = 5500
$ws.Cells.Item(6, 13).Value2 = 5250
$ws.Cells.Item(6, 14).Value2 = "$/docena de matas"
$ws.Cells.Item(6, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(6, 16).Value2 = 875
$ws.Cells.Item(6, 17).Value2 = 6
$ws.Cells.Item(6, 18).Value2 = "Hortaliza"

# Make sure the new date cell keeps the same date number format as the
# rest of column D (style copied automatically by Insert, but set it
# explicitly too in case it wasn't carried over).
$ws.Cells.Item(6, 4).NumberFormat = $ws.Cells.Item(7, 4).NumberFormat

Write-Output "Done. New row count: $($ws.UsedRange.Rows.Count)"
